# Generate Report for Handback
#
# This script reflects that the 882b804b-... and af6f3f4b-... files have now
# been handed back (for both the zh-cn and de-de locales). It updates:
#  - the Overview sheet's per-locale status/date columns for those two files
#  - each locale sheet's Status column, Latest Target File, Latest Handback
#    File and Latest Handback DateTime columns for those two rows

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (882b804b) and 5 (af6f3f4b) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn (E) and de-de (F) columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $statusHandedBack
$wsOverview.Range("F4").Value = $statusHandedBack
$wsOverview.Range("E5").Value = $statusHandedBack
$wsOverview.Range("F5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (882b804b) and 5 (af6f3f4b)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusHandedBack
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d5d5f6d40eb297e32f0b77ade11a65d13e046026/e2e/882b804b-fd0d-4f7b-b456-97acbbdf8fb1.md", "", "", "882b804b-fd0d-4f7b-b456-97acbbdf8fb1.md") | Out-Null
$wsZhCn.Range("J4").Value = "882b804b-fd0d-4f7b-b456-97acbbdf8fb1.fdca27ef944a4b48451a6e11bbd0813e9d93073a.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-21 06:38:28"

$wsZhCn.Range("C5").Value = $statusHandedBack
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d5d5f6d40eb297e32f0b77ade11a65d13e046026/e2e/af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.md", "", "", "af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.md") | Out-Null
$wsZhCn.Range("J5").Value = "af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.fec5cc13855640bb598568739fe5e50096a2629f.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-21 06:38:28"

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (882b804b) and 5 (af6f3f4b)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusHandedBack
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5cf30bc51e139fad342e8acb48fb6d7740e38af9/e2e/882b804b-fd0d-4f7b-b456-97acbbdf8fb1.md", "", "", "882b804b-fd0d-4f7b-b456-97acbbdf8fb1.md") | Out-Null
$wsDeDe.Range("J4").Value = "882b804b-fd0d-4f7b-b456-97acbbdf8fb1.fdca27ef944a4b48451a6e11bbd0813e9d93073a.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-21 06:38:35"

$wsDeDe.Range("C5").Value = $statusHandedBack
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5cf30bc51e139fad342e8acb48fb6d7740e38af9/e2e/af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.md", "", "", "af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.md") | Out-Null
$wsDeDe.Range("J5").Value = "af6f3f4b-43a6-4e1f-b4bf-89198fd93d20.fec5cc13855640bb598568739fe5e50096a2629f.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-21 06:38:35"
